# The post previously stored at row 239 ("自分を信じなさい。人の言うことは尽きないのだから")
# was removed after publishing. Delete that entire row so every following
# row (240-300) shifts up by one, matching the renumbering seen in the diff,
# and the sheet's used range shrinks from A1:C300 to A1:C299.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(239).Delete()
